$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 524.4
$ws.Range("I43").Value = 472
$ws.Range("J43").Value = 559.3333
$ws.Range("K43").Value = 472
$ws.Range("L43").Value = 559.3333
$ws.Range("M43").Value = -403
$ws.Range("N43").Value = -697.3333
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H48").Value = 5271.4287
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 5271.4287
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 15814.2861
$ws.Range("N48").Value = -16398.2861
$ws.Range("H51").Value = 1685.5714
$ws.Range("J51").Value = 1560
$ws.Range("L51").Value = 1560
$ws.Range("N51").Value = -2528
$ws.Range("H56").Value = 5271.4287
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 5271.4287
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 15814.2861
$ws.Range("N56").Value = -16882.2861
$ws.Range("H98").Value = 1622.8182
$ws.Range("I98").Value = 1251.9375
$ws.Range("J98").Value = 2611.8333
$ws.Range("K98").Value = 1251.9375
$ws.Range("L98").Value = 2611.8333
$ws.Range("M98").Value = 246.0625
$ws.Range("N98").Value = -5607.8333
$ws.Range("H122").Value = 1622.8182
$ws.Range("I122").Value = 1251.9375
$ws.Range("J122").Value = 2611.8333
$ws.Range("K122").Value = 3755.8125
$ws.Range("L122").Value = 7835.499899999999
$ws.Range("M122").Value = -1305.8125
$ws.Range("N122").Value = -12735.4999
$ws.Range("H131").Value = 28703.244
$ws.Range("I131").Value = 33206.453
$ws.Range("J131").Value = 5436.6665
$ws.Range("K131").Value = 99619.359
$ws.Range("L131").Value = 16309.9995
$ws.Range("M131").Value = -94579.359
$ws.Range("N131").Value = -26389.9995
$ws.Range("H135").Value = 955
$ws.Range("I135").Value = 679.75
$ws.Range("J135").Value = 1780.75
$ws.Range("K135").Value = 6117.75
$ws.Range("L135").Value = 16026.75
$ws.Range("M135").Value = -3582.75
$ws.Range("N135").Value = -21096.75
$ws.Range("H137").Value = 33198.656
$ws.Range("I137").Value = 1666.5264
$ws.Range("K137").Value = 4999.5792
$ws.Range("M137").Value = -2449.5792
$ws.Range("M45").ClearContents()
$ws.Range("M48").ClearContents()
$ws.Range("M56").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 9445
$ws.Range("I8").Value = 2000
$ws.Range("J8").Value = 16890
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 16890
$ws.Range("M8").Value = -1856
$ws.Range("N8").Value = -17178
$ws.Range("H13").Value = 1505001.5
$ws.Range("I13").Value = 3000003
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 3000003
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = -2999859
$ws.Range("N13").Value = -10288
$ws.Range("H74").Value = 8229.866
$ws.Range("I74").Value = 9188.385
$ws.Range("J74").Value = 1999.5
$ws.Range("K74").Value = 9188.385
$ws.Range("L74").Value = 1999.5
$ws.Range("M74").Value = -8314.385
$ws.Range("N74").Value = -3747.5
$ws.Range("H77").Value = 8229.866
$ws.Range("I77").Value = 9188.385
$ws.Range("J77").Value = 1999.5
$ws.Range("K77").Value = 45941.925
$ws.Range("L77").Value = 9997.5
$ws.Range("M77").Value = -41573.925
$ws.Range("N77").Value = -18733.5
$ws.Range("H110").Value = 1967.8889
$ws.Range("I110").Value = 2203.6667
$ws.Range("J110").Value = 1850
$ws.Range("K110").Value = 2203.6667
$ws.Range("L110").Value = 1850
$ws.Range("M110").Value = -158.6667000000002
$ws.Range("N110").Value = -5940
$ws.Range("H132").Value = 4628.25
$ws.Range("I132").Value = 4602.4
$ws.Range("J132").Value = 4671.3335
$ws.Range("K132").Value = 13807.2
$ws.Range("L132").Value = 14014.0005
$ws.Range("M132").Value = -11277.2
$ws.Range("N132").Value = -19074.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 227.35715
$ws.Range("J94").Value = 220.625
$ws.Range("L94").Value = 220.625
$ws.Range("N94").Value = -1122.625
$ws.Range("H107").Value = 2121.8125
$ws.Range("I107").Value = 2722.2222
$ws.Range("J107").Value = 1349.8572
$ws.Range("K107").Value = 2722.2222
$ws.Range("L107").Value = 1349.8572
$ws.Range("M107").Value = -802.2222000000002
$ws.Range("N107").Value = -5189.8572
$ws.Range("H134").Value = 2063.0588
$ws.Range("I134").Value = 2063.0588
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6189.176399999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3654.176399999999
$ws.Range("N134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1133.3636
$ws.Range("I22").Value = 994.7143
$ws.Range("J22").Value = 1376
$ws.Range("K22").Value = 994.7143
$ws.Range("L22").Value = 1376
$ws.Range("M22").Value = -644.7143
$ws.Range("N22").Value = -2076
$ws.Range("H31").Value = 2731
$ws.Range("I31").Value = 2051.919
$ws.Range("J31").Value = 5243.6
$ws.Range("K31").Value = 2051.919
$ws.Range("L31").Value = 5243.6
$ws.Range("M31").Value = -1756.919
$ws.Range("N31").Value = -5833.6
$ws.Range("H34").Value = 2731
$ws.Range("I34").Value = 2051.919
$ws.Range("J34").Value = 5243.6
$ws.Range("K34").Value = 2051.919
$ws.Range("L34").Value = 5243.6
$ws.Range("M34").Value = -1849.919
$ws.Range("N34").Value = -5647.6
$ws.Range("H134").Value = 3603.2307
$ws.Range("I134").Value = 1802.5454
$ws.Range("J134").Value = 13507
$ws.Range("K134").Value = 5407.6362
$ws.Range("L134").Value = 40521
$ws.Range("M134").Value = -2872.6362
$ws.Range("N134").Value = -45591

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1472.8182
$ws.Range("I86").Value = 802
$ws.Range("J86").Value = 1539.9
$ws.Range("K86").Value = 2406
$ws.Range("L86").Value = 4619.700000000001
$ws.Range("M86").Value = -1220
$ws.Range("N86").Value = -6991.700000000001
$ws.Range("H89").Value = 1472.8182
$ws.Range("I89").Value = 802
$ws.Range("J89").Value = 1539.9
$ws.Range("K89").Value = 7218
$ws.Range("L89").Value = 13859.1
$ws.Range("M89").Value = -1290
$ws.Range("N89").Value = -25715.1
$ws.Range("H98").Value = 252.125
$ws.Range("I98").Value = 255.8
$ws.Range("J98").Value = 246
$ws.Range("K98").Value = 767.4000000000001
$ws.Range("L98").Value = 738
$ws.Range("M98").Value = 730.5999999999999
$ws.Range("N98").Value = -3734
$ws.Range("H132").Value = 1123523.6
$ws.Range("I132").Value = 1224.5
$ws.Range("J132").Value = 2021363
$ws.Range("K132").Value = 11020.5
$ws.Range("L132").Value = 18192267
$ws.Range("M132").Value = -8490.5
$ws.Range("N132").Value = -18197327

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.11539
$ws.Range("I2").Value = 56.1
$ws.Range("J2").Value = 139
$ws.Range("K2").Value = 56.1
$ws.Range("L2").Value = 139
$ws.Range("M2").Value = 56.9
$ws.Range("N2").Value = -365
$ws.Range("H102").Value = 2072.5881
$ws.Range("I102").Value = 1902.7858
$ws.Range("J102").Value = 2865
$ws.Range("K102").Value = 1902.7858
$ws.Range("L102").Value = 2865
$ws.Range("M102").Value = -280.7858000000001
$ws.Range("N102").Value = -6109
$ws.Range("H107").Value = 166.70589
$ws.Range("I107").Value = 181
$ws.Range("K107").Value = 181
$ws.Range("M107").Value = 1739
$ws.Range("H113").Value = 11336.917
$ws.Range("I113").Value = 2227.111
$ws.Range("J113").Value = 38666.332
$ws.Range("K113").Value = 2227.111
$ws.Range("L113").Value = 38666.332
$ws.Range("M113").Value = -57.11099999999988
$ws.Range("N113").Value = -43006.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4078.35
$ws.Range("I40").Value = 3568.647
$ws.Range("J40").Value = 6966.6665
$ws.Range("K40").Value = 3568.647
$ws.Range("L40").Value = 6966.6665
$ws.Range("M40").Value = -3432.647
$ws.Range("N40").Value = -7238.6665
$ws.Range("H61").Value = 11682.223
$ws.Range("I61").Value = 33663.332
$ws.Range("J61").Value = 691.6667
$ws.Range("K61").Value = 33663.332
$ws.Range("L61").Value = 691.6667
$ws.Range("M61").Value = -33461.332
$ws.Range("N61").Value = -1095.6667
$ws.Range("H113").Value = 11682.223
$ws.Range("I113").Value = 33663.332
$ws.Range("J113").Value = 691.6667
$ws.Range("K113").Value = 33663.332
$ws.Range("L113").Value = 691.6667
$ws.Range("M113").Value = -31493.332
$ws.Range("N113").Value = -5031.6667
$ws.Range("H132").Value = 7152.2666
$ws.Range("I132").Value = 6216.727
$ws.Range("K132").Value = 18650.181
$ws.Range("M132").Value = -16120.181

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 435.66666
$ws.Range("I113").Value = 428
$ws.Range("K113").Value = 1284
$ws.Range("M113").Value = 886
$ws.Range("H122").Value = 1987.0667
$ws.Range("I122").Value = 1604.8422
$ws.Range("J122").Value = 2647.2727
$ws.Range("K122").Value = 4814.5266
$ws.Range("L122").Value = 7941.8181
$ws.Range("M122").Value = -2364.5266
$ws.Range("N122").Value = -12841.8181
